$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.984.40"
$ws.Range("E2").Value = "  +0.10%  "

$ws.Range("D3").Value = "2.238.72"
$ws.Range("E3").Value = "  -0.19%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "306.75"
$ws.Range("E5").Value = "  -3.89%  "

$ws.Range("D6").Value = "94.95"
$ws.Range("E6").Value = "  -5.12%  "

$ws.Range("E7").Value = "  -0.41%  "

$ws.Range("E8").Value = "  +0.18%  "

$ws.Range("D9").Value = "0.522"
$ws.Range("E9").Value = "  -4.02%  "

$ws.Range("D10").Value = "34.81"
$ws.Range("E10").Value = "  -4.75%  "

$ws.Range("D11").Value = "0.0809"
$ws.Range("E11").Value = "  -1.77%  "

$ws.Range("E12").Value = "  -3.41%  "

$ws.Range("E13").Value = "  -1.24%  "

$ws.Range("D14").Value = "2.579.61"
$ws.Range("E14").Value = "  -0.25%  "

$ws.Range("D15").Value = "2.310.35"
$ws.Range("E15").Value = "  +2.74%  "

$ws.Range("E16").Value = "  -2.42%  "

$ws.Range("E17").Value = "  -4.36%  "

$ws.Range("D18").Value = "43.878.60"
$ws.Range("E18").Value = "  +0.05%  "

$ws.Range("D19").Value = "0.0₃0960"
$ws.Range("E19").Value = "  -1.21%  "

$ws.Range("D20").Value = "12.11"
$ws.Range("E20").Value = "  -8.88%  "

$ws.Range("D21").Value = "6.27"
$ws.Range("E21").Value = "  -1.99%  "

$ws.Range("D22").Value = "65.11"
$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("D23").Value = "236.77"
$ws.Range("E23").Value = "  +1.65%  "

$ws.Range("E24").Value = "  -4.74%  "

$ws.Range("E25").Value = "  -4.94%  "

$ws.Range("E26").Value = "  -0.13%  "

$ws.Range("E27").Value = "  -5.37%  "

$ws.Range("E28").Value = "  -0.13%  "

$ws.Range("E29").Value = "  -2.12%  "

$ws.Range("D30").Value = "6.04"
$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("D31").Value = "19.92"
$ws.Range("E31").Value = "  -0.42%  "

$ws.Range("D32").Value = "152.87"
$ws.Range("E32").Value = "  -3.41%  "

$ws.Range("D33").Value = "0.0803"
$ws.Range("E33").Value = "  -4.42%  "

$ws.Range("D34").Value = "3.25"
$ws.Range("E34").Value = "  +4.23%  "

$ws.Range("D35").Value = "2.58"
$ws.Range("E35").Value = "  -3.76%  "

$ws.Range("E36").Value = "  -3.80%  "

$ws.Range("D37").Value = "0.119"
$ws.Range("E37").Value = "  +0.77%  "

$ws.Range("E38").Value = "  -7.17%  "

$ws.Range("D39").Value = "15.11"
$ws.Range("E39").Value = "  -6.75%  "

$ws.Range("D40").Value = "3.86"
$ws.Range("E40").Value = "  -6.78%  "

$ws.Range("E41").Value = "  -8.33%  "

$ws.Range("E42").Value = "  -3.71%  "

$ws.Range("E43").Value = "  +0.24%  "

$ws.Range("D44").Value = "1.724.60"
$ws.Range("E44").Value = "  -1.91%  "

$ws.Range("D45").Value = "85.42"
$ws.Range("E45").Value = "  +5.89%  "

$ws.Range("E46").Value = "  -3.34%  "

$ws.Range("D47").Value = "100.14"
$ws.Range("E47").Value = "  -2.58%  "

$ws.Range("E48").Value = "  -3.72%  "

$ws.Range("D49").Value = "69.14"
$ws.Range("E49").Value = "  -6.65%  "

$ws.Range("D50").Value = "8.08"
$ws.Range("E50").Value = "  -2.04%  "

$ws.Range("D51").Value = "54.27"
$ws.Range("E51").Value = "  -4.46%  "
